# Fix du téléchargement des toutes les sessions
# Rebuilds the "Campagne" export sheet: proper header row (Token / Coords /
# Association Ligne / Association Col / Moyenne / Temps total / Choix Final /
# Id Campagne) with bold + bordered + centered/top-aligned styling, and the
# per-token result rows collected for this session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old single-row sample content first so none of its stale
# shared strings (old Token/Coords/... placeholder text) linger in the
# rebuilt table.
$ws.Range("A1:G1").ClearContents()

# ---------------------------------------------------------------------
# Header row (row 1) values
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 1).Value = "Token"
$ws.Cells.Item(1, 2).Value = "Coords"
$ws.Cells.Item(1, 3).Value = "Association Ligne"
$ws.Cells.Item(1, 4).Value = "Association Col"
$ws.Cells.Item(1, 5).Value = "Moyenne"
$ws.Cells.Item(1, 6).Value = "Temps total"
$ws.Cells.Item(1, 7).Value = "Choix Final"
$ws.Cells.Item(1, 8).Value = "Id Campagne"

# Build the header look (bold, thin box border, centered + top-aligned) on
# a single cell, then fan it out with a format-only paste so every header
# cell shares the exact same style entry instead of each property mutation
# minting its own cellXf.
$headerTemplate = $ws.Cells.Item(1, 1)
$headerTemplate.Font.Bold = $true
$headerTemplate.Borders.LineStyle = 1
$headerTemplate.HorizontalAlignment = -4108
$headerTemplate.VerticalAlignment = -4160
$headerTemplate.Copy()
$ws.Range("A1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Row 2 - token C1ID431991
# ---------------------------------------------------------------------
$ws.Cells.Item(2, 1).Value = "C1ID431991"
$ws.Cells.Item(2, 2).Value = "2:1', 3:2', 3:3', 4:3', 4:4', 3:5'"
$ws.Cells.Item(2, 3).Value = "Shining', Blues brothers', Hellraiser', Hellraiser', 12 hommes en coleres', Massacre \xe0 la tron\xe7onneuse'"
$ws.Cells.Item(2, 4).Value = "R\xe9alisateur', Ann\xe9e', Ann\xe9e', Genre', Genre', Ann\xe9e'"
$ws.Cells.Item(2, 5).Value = "Massacre à la tronçonneuse"
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 3771

# ---------------------------------------------------------------------
# Row 3 - token C1IDe1f63e
# ---------------------------------------------------------------------
$ws.Cells.Item(3, 1).Value = "C1IDe1f63e"
$ws.Cells.Item(3, 2).Value = "2:1', 2:2', 3:3', 3:4', 1:5'"
$ws.Cells.Item(3, 3).Value = "Hellraiser', Blues brothers', Massacre \xe0 la tron\xe7onneuse', Shining', 12 hommes en coleres'"
$ws.Cells.Item(3, 4).Value = "R\xe9alisateur', R\xe9alisateur', Ann\xe9e', Ann\xe9e', Affiche'"
$ws.Cells.Item(3, 5).Value = "12 hommes en coleres"
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 3784

# ---------------------------------------------------------------------
# Row 4 - token C1ID7f262d
# ---------------------------------------------------------------------
$ws.Cells.Item(4, 1).Value = "C1ID7f262d"
$ws.Cells.Item(4, 2).Value = "2:1', 2:2', 2:3', 2:4', 2:5'"
$ws.Cells.Item(4, 3).Value = "Massacre \xe0 la tron\xe7onneuse', Shining', 12 hommes en coleres', Hellraiser', Blues brothers'"
$ws.Cells.Item(4, 4).Value = "Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e'"
$ws.Cells.Item(4, 5).Value = "Blues brothers"
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 2827

# ---------------------------------------------------------------------
# Rows 5-6 - duplicate "Association" echoes left over from the grid
# (same Coords/Association text as row 4, only the running total column
# is populated besides that)
# ---------------------------------------------------------------------
$ws.Cells.Item(5, 2).Value = "2:1', 2:2', 2:3', 2:4', 2:5'"
$ws.Cells.Item(5, 3).Value = "Massacre \xe0 la tron\xe7onneuse', Shining', 12 hommes en coleres', Hellraiser', Blues brothers'"
$ws.Cells.Item(5, 4).Value = "Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e'"
$ws.Cells.Item(5, 6).Value = 1

$ws.Cells.Item(6, 2).Value = "2:1', 2:2', 2:3', 2:4', 2:5'"
$ws.Cells.Item(6, 3).Value = "Massacre \xe0 la tron\xe7onneuse', Shining', 12 hommes en coleres', Hellraiser', Blues brothers'"
$ws.Cells.Item(6, 4).Value = "Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e', Ann\xe9e'"
$ws.Cells.Item(6, 6).Value = 1

# ---------------------------------------------------------------------
# Touch the remaining blank cells of rows 2-9 (no value, no style change)
# so they stay present in the sheet grid, matching the template's
# pre-sized A1:H9 range.
# ---------------------------------------------------------------------
$blankCells = @(
    @(2, 8), @(3, 8), @(4, 8),
    @(5, 1), @(5, 5), @(5, 7), @(5, 8),
    @(6, 1), @(6, 5), @(6, 7), @(6, 8),
    @(7, 1), @(7, 2), @(7, 3), @(7, 4), @(7, 5), @(7, 6), @(7, 7), @(7, 8),
    @(8, 1), @(8, 2), @(8, 3), @(8, 4), @(8, 5), @(8, 6), @(8, 7), @(8, 8),
    @(9, 1), @(9, 2), @(9, 3), @(9, 4), @(9, 5), @(9, 6), @(9, 7), @(9, 8)
)
foreach ($pos in $blankCells) {
    $ws.Cells.Item($pos[0], $pos[1]).Font.Bold = $false
}
